# issue #5: stock data output to json file
#
# The "股票" (stock) worksheet gains a new "property_category" column
# (constant value "stock") inserted right before the existing "date"
# column. The date / legislator_name / legislator_id columns shift one
# column to the right as a result.
#
# Three company-name strings in the stock sheet also had a stray space
# removed as part of this change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- fix stray-space typos in company names (column B, name) ---------
$ws.Range("B2").Value = "台灣卜蜂企業股份有限公司"
$ws.Range("B5").Value = "同泰電子科技股份有限公司"
$ws.Range("B6").Value = "上福全球科技股份有限公司"

# --- insert the new property_category column before the date column --
$ws.Columns("H").Insert()

$ws.Range("H1").Value = "property_category"
$ws.Range("H2:H7").Value = "stock"
